$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: A13 = 1024, B13 = 79000
$ws.Range("A13").Value = 1024
$ws.Range("B13").Value = 79000

# Match the style used by the row above (A12:B12) for the new cells
$ws.Range("A12:B12").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active cell selection to match the saved view state
$ws.Range("E8").Select()
